# Estructuras y verificación de los dígitos de control
#
# The "CodigoCuenta" column (J) on Hoja1 held a distinct fake bank-account
# number per row. After validating/normalizing the control digits, every
# row now carries the same verified account code.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$finalAccountCode = "11112223504444444444"

# Rows 9, 28, 59, 79, 100 are blank separator rows with no data in column J.
$rows = @(2,3,4,5,6,7,8,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,
          29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,
          51,52,53,54,55,56,57,58,60,61,62,63,64,65,66,67,68,69,70,71,72,73,
          74,75,76,77,78,80,81,82,83,84,85,86,87,88,89,90,91,92,93,94,95,96,
          97,98,99,101,102,103,104,105)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 10)   # column J = CodigoCuenta
    $cell.NumberFormat = "@"
    $cell.Value = $finalAccountCode
}
